# IS_training_format.xlsx -- simplify the training sheet down to a plain
# 4-column table (sl no. / training name / description / Remarks) and drop
# the old wide, merged, bordered "multiple dept" layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop both header rows (the two-row merged header with the duplicated
# "Grade" column J/K was the source of the "multiple dept" bug) and the
# extra columns E:O that aren't needed any more.
$ws.Rows("1:2").Delete()
$ws.Columns("E:O").Delete()

# Write the new, simple header row.
$ws.Range("A1").Value = "sl no."
$ws.Range("B1").Value = "training name"
$ws.Range("C1").Value = "description "
$ws.Range("D1").Value = "Remarks"

# Size the new columns to fit the new header text.
$ws.Columns("B:B").ColumnWidth = 15.166666666666666
$ws.Columns("C:C").ColumnWidth = 12.276041666666666
$ws.Columns("D:D").ColumnWidth = 20.608072916666668

# Match the saved selection/cursor position.
$ws.Range("D2").Select() | Out-Null
